# Apply review-db update: two rows' "confirm" flag flips from yes->no,
# and two brand-new review rows are appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10 and 11: column G ("blue"/confirm) changes from "yes" to "no"
$ws.Range("G10").Value = "no"
$ws.Range("G11").Value = "no"

# New row 23
$ws.Range("A23").Value = "com.hamxa.shaynachim"
$ws.Range("A23").Font.Name = "Mangal"
$ws.Range("A23").Font.Size = 10
$ws.Range("B23").Value = "bitcoin"
$ws.Range("C23").Value = "itaisenior@gmail.com"
$ws.Range("D23").Value = "vikicrestina@gmail.com"
$ws.Range("E23").Value = "27/5/2019 15:59"
$ws.Range("F23").Value = "If you really want to know what is bitcoin and blockchain use this app"
$ws.Range("G23").Value = "no"

# New row 24
$ws.Range("A24").Value = "com.hamxa.shaynachim"
$ws.Range("A24").Font.Name = "Mangal"
$ws.Range("A24").Font.Size = 10
$ws.Range("B24").Value = "bitcoin"
$ws.Range("C24").Value = "leviadlevi22@gmail.com"
$ws.Range("D24").Value = "gazittalia1@gmail.com"
$ws.Range("E24").Value = "27/5/2019 15:59"
$ws.Range("F24").Value = "awesome content and great written. Exactly in the spot"
$ws.Range("G24").Value = "no"

# Hyperlinks on the email/recovery columns for the two new rows (mailto:), matching
# the existing pattern used throughout the sheet. Adding a hyperlink repaints the
# cell with the blue underlined "Hyperlink" style, so restore the plain look used
# by every other linked cell in this sheet (Calibri 11 black, centered).
$ws.Hyperlinks.Add($ws.Range("C23"), "mailto:itaisenior@gmail.com", "", "", "itaisenior@gmail.com")
$ws.Range("C23").Font.Name = "Calibri"
$ws.Range("C23").Font.Size = 11
$ws.Range("C23").Font.Underline = $false
$ws.Range("C23").Font.Color = 0
$ws.Range("C23").HorizontalAlignment = -4108

$ws.Hyperlinks.Add($ws.Range("D23"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com")
$ws.Range("D23").Font.Name = "Calibri"
$ws.Range("D23").Font.Size = 11
$ws.Range("D23").Font.Underline = $false
$ws.Range("D23").Font.Color = 0
$ws.Range("D23").HorizontalAlignment = -4108

$ws.Hyperlinks.Add($ws.Range("C24"), "mailto:leviadlevi22@gmail.com", "", "", "leviadlevi22@gmail.com")
$ws.Range("C24").Font.Name = "Calibri"
$ws.Range("C24").Font.Size = 11
$ws.Range("C24").Font.Underline = $false
$ws.Range("C24").Font.Color = 0
$ws.Range("C24").HorizontalAlignment = -4108

$ws.Hyperlinks.Add($ws.Range("D24"), "mailto:gazittalia1@gmail.com", "", "", "gazittalia1@gmail.com")
$ws.Range("D24").Font.Name = "Calibri"
$ws.Range("D24").Font.Size = 11
$ws.Range("D24").Font.Underline = $false
$ws.Range("D24").Font.Color = 0
$ws.Range("D24").HorizontalAlignment = -4108

# Leave the cursor where the author's session ended up
$ws.Range("F25").Select()
